$wb = $excel.ActiveWorkbook

# Tabelle2 ("cats" translation row) gets its description strings touched up:
# punctuation added to the "Dont use" / "this row" columns for that entry.
$ws = $wb.Worksheets.Item("Tabelle2")
$ws.Range("D3").Value = "Katzenbeschreibung!"
$ws.Range("E3").Value = "Noch ne Katzenbeschreibung…"

# Reflect the author's new selection/window state when they saved the file.
$ws.Activate()
$ws.Range("E3").Select()
$excel.ActiveWindow.WindowState = [Microsoft.Office.Interop.Excel.XlWindowState]::xlMaximized
